$d = $word.ActiveDocument

function Replace-UniqueText($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Date heading
Replace-UniqueText "2024-08-13 Tuesday" "2024-08-14 Wednesday"

# Row 1 of the table contains two identical cells ("44÷9=4, 8") that must
# receive different replacement values, so address them positionally via
# the table cell ranges instead of a document-wide Find/Replace.
$tbl = $d.Tables(1)

$cell3 = $tbl.Cell(1, 3)
$r3 = $d.Range($cell3.Range.Start, $cell3.Range.End)
$r3.Find.Execute("44÷9=4, 8", $true, $false, $false, $false, $false, $true, 0, $false, "29÷5=5, 4", 2) | Out-Null

$tbl = $d.Tables(1)
$cell4 = $tbl.Cell(1, 4)
$r4 = $d.Range($cell4.Range.Start, $cell4.Range.End)
$r4.Find.Execute("44÷9=4, 8", $true, $false, $false, $false, $false, $true, 0, $false, "57÷7=8, 1", 2) | Out-Null

# Remaining unique division problems
Replace-UniqueText "34÷6=5, 4" "51÷8=6, 3"
Replace-UniqueText "42÷6=7, 0" "50÷3=16, 2"
Replace-UniqueText "45÷9=5, 0" "76÷8=9, 4"
Replace-UniqueText "71÷4=17, 3" "67÷9=7, 4"
Replace-UniqueText "18÷6=3, 0" "78÷6=13, 0"
Replace-UniqueText "62÷2=31, 0" "74÷2=37, 0"
Replace-UniqueText "43÷8=5, 3" "59÷7=8, 3"
Replace-UniqueText "50÷4=12, 2" "55÷9=6, 1"
Replace-UniqueText "52÷7=7, 3" "77÷5=15, 2"
Replace-UniqueText "42÷9=4, 6" "36÷7=5, 1"
Replace-UniqueText "14÷9=1, 5" "84÷4=21, 0"
Replace-UniqueText "67÷8=8, 3" "68÷2=34, 0"
Replace-UniqueText "70÷6=11, 4" "26÷8=3, 2"
Replace-UniqueText "88÷7=12, 4" "48÷3=16, 0"
Replace-UniqueText "73÷8=9, 1" "21÷5=4, 1"
Replace-UniqueText "64÷7=9, 1" "84÷6=14, 0"
Replace-UniqueText "44÷8=5, 4" "62÷8=7, 6"
Replace-UniqueText "31÷2=15, 1" "61÷3=20, 1"
Replace-UniqueText "31÷3=10, 1" "25÷9=2, 7"
Replace-UniqueText "49÷2=24, 1" "22÷2=11, 0"
Replace-UniqueText "47÷3=15, 2" "49÷4=12, 1"
Replace-UniqueText "45÷6=7, 3" "31÷6=5, 1"
Replace-UniqueText "76÷2=38, 0" "58÷6=9, 4"

Write-Output "All replacements applied"
